$d = $word.ActiveDocument

# Collect the indices (1-based) of every paragraph whose text is exactly
# "Rex change it 2" (the paragraphs built from the two runs
# "Rex " + "change it 2"). Paragraph.Range.Text includes the trailing
# paragraph-mark character, so strip that before comparing.
$targets = @()
$idx = 1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    $t = $t.TrimEnd([char]13, [char]7)
    if ($t -eq "Rex change it 2") {
        $targets += , $idx
    }
    $idx = $idx + 1
}

# Process from the last match to the first so that inserting new
# paragraphs above earlier matches never invalidates indices we still
# need to visit.
if ($targets.Count -gt 0) {
    $targets = $targets[($targets.Count - 1)..0]
}

foreach ($i in $targets) {
    $target = $d.Paragraphs($i)

    # 1) Insert a brand new paragraph directly above the matched one,
    #    containing the plain text "Rex change it 2".
    $target.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs($i)
    $newPara.Range.Text = "Rex change it 2"

    # 2) In the original paragraph (now one position further down),
    #    change "change it 2" to "change it 3".
    $origPara = $d.Paragraphs($i + 1)
    $replaceRange = $origPara.Range.Duplicate
    $replaceRange.Find.Execute("change it 2", $false, $false, $false, $false, $false, `
                                $true, 1, $false, "change it 3", 2)
}
